$wb = $excel.ActiveWorkbook

# Rename the "read me" sheet to "ReadMe"
$ws = $wb.ActiveSheet
$ws.Name = "ReadMe"

# Update row 1
$ws.Range("A1").Value = "Data Dictionary"
$ws.Range("B1").Value = "...2"
$ws.Range("C1").Value = "...3"
$ws.Range("D1").Value = "...4"

# Row 2
$ws.Range("A2").Value = "Authors:"
$ws.Range("C2").Value = "R Sapir-Pichhadze, E Gitelman, S El Bouzaidi Tiali, G Fabre, J Laforme"

# Row 4
$ws.Range("A4").Value = "Aim of the document:"
$ws.Range("C4").Value = "Provide a standard metadata model for the input dataset in the banffIT package. It describes the format and the constraints the input dataset must follow to be able to run the diagnosis assignement process"

# Row 6
$ws.Range("A6").Value = "Current version:"
$ws.Range("C6").Value = "Version 1.0"

# Row 7
$ws.Range("A7").Value = "Banff Classification version:"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "2017"
$ws.Range("C7").ClearFormats()

# Row 9
$ws.Range("A9").Value = "Spreadsheets:"
$ws.Range("C9").Value = "Description"

# Row 10
$ws.Range("B10").Value = "ReadMe"
$ws.Range("C10").Value = "Description of the document aims, version, and content"

# Row 11
$ws.Range("B11").Value = "Variables"
$ws.Range("C11").Value = "Metadata of the list of minimum input variables to be present in the input dataset along with the contraints they must follow in order to be able to assign diagnoses of the Banff Classification using the banffIT package. It also contains the metadata of all possible output variables."

# Row 12
$ws.Range("B12").Value = "Categories"
$ws.Range("C12").Value = "Code, label, and description of each categorical variables that should be present in the input dataset and that can be present in the output dataset"

# Row 14
$ws.Range("A14").Value = "Version history:"
$ws.Range("C14").Value = "Description"
$ws.Range("D14").Value = "Date"

# Clear out old row 15/16 content first, then set new row 15
$ws.Range("A15").Value = ""
$ws.Range("A16").Value = ""
$ws.Range("B16").Value = ""

# Row 15 (new content)
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "1.0"
$ws.Range("B15").ClearFormats()
$ws.Range("C15").Value = "Initial version of the data dictionary"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "45415"
$ws.Range("D15").ClearFormats()
